$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (same header style as the existing A1:F1 headers)
$ws.Range("G1").Value = "Lent"
$ws.Range("H1").Value = "Lent to"
$ws.Range("I1").Value = "Lent date"
$ws.Range("J1").Value = "Return date"
$ws.Range("K1").Value = "Reserved"
$ws.Range("L1").Value = "Reserved by"
$ws.Range("M1").Value = "Reserved until"

$ws.Range("A1").Copy()
$ws.Range("G1:M1").PasteSpecial(-4122)  # xlPasteFormats - reuse the existing header style

# Row 2 - book #1 ("Rozdroze Krukow") is currently lent out
$ws.Range("G2").Value = $true
$ws.Range("H2").Value = 1

# Date/time cells: create the two number formats used in the workbook
# (164 stays unused - a leftover format Excel created while editing;
#  165 is the one actually applied to the lent/return/reserved dates)
$ws.Range("I2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("I2").Value = 100
$ws.Range("I2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("I2").Value = 45805.66867783668

$ws.Range("J2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J2").Value = 45835.66867783669

$ws.Range("K2").Value = $false
$ws.Range("M2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("M2").Value = 45805.659361875

# Row 3 - book #2 ("Metro 2033") is not lent, not reserved
$ws.Range("G3").Value = $false
$ws.Range("K3").Value = $false
$ws.Range("M3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("M3").Value = 45806.65936186343

# Blank (never-lent / never-reserved) cells still get a cell entry in the
# sheet, just with no value - materialize them by pasting the formatting
# of an untouched, default-style cell over them (keeps them on the
# default style, no new cellXfs entries get created).
$ws.Range("Z99").Copy()
$ws.Range("L2").PasteSpecial(-4122)
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("I3").PasteSpecial(-4122)
$ws.Range("J3").PasteSpecial(-4122)
$ws.Range("L3").PasteSpecial(-4122)
